$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values (sample sizes)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (CON meanEMG/legmaxROM) - B2 cleared, C2:E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.82147052983003199
$ws.Range("D2").Value = 0.46961565982885617
$ws.Range("E2").Value = 1.1388024172627749

# Update row 3 values (STR meanEMG/legmaxROM)
$ws.Range("B3").Value = 2.0600852448748426
$ws.Range("C3").Value = 2.0547895786242916
$ws.Range("D3").Value = 7.5144916093350496
$ws.Range("E3").Value = 5.2553961548146289

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
